$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.317.92"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.876.01"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'0.7106"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "'242.87"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.07991"
$ws.Range("E8").Value = "  +2.51%  "
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").Value = "'24.95"
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("D11").Value = "'0.08240"
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("D12").Value = "1.896.12"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("D13").Value = "'5.241"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").Value = "'94.54"
$ws.Range("E14").Value = "  +3.73%  "
$ws.Range("D15").Value = "'0.7116"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "'6.349"
$ws.Range("E16").Value = "  +4.04%  "
$ws.Range("D17").Value = "'0.000008540"
$ws.Range("E17").Value = "  +3.97%  "
$ws.Range("D18").Value = "29.340.01"
$ws.Range("D19").Value = "'244.82"
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("D20").Value = "2.153.63"
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'7.779"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "'0.1552"
$ws.Range("E25").Value = "  -2.95%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'162.59"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'9.043"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "'18.50"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'4.413"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("D31").Value = "'4.309"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'1.183"
$ws.Range("E32").Value = "  -8.43%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.05380"
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("D34").Value = "'1.935"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "'0.7651"
$ws.Range("E35").Value = "  +2.69%  "
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("D37").Value = "'2.688"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("D38").Value = "'0.01879"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("D39").Value = "1.258.65"
$ws.Range("E39").Value = "  +2.67%  "
$ws.Range("D40").Value = "'2.749"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").Value = "'6.498"
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("D42").Value = "'0.9184"
$ws.Range("E42").Value = "  +3.40%  "
$ws.Range("D43").Value = "'112.91"
$ws.Range("E43").Value = "  +1.80%  "
$ws.Range("D44").Value = "'74.17"
$ws.Range("E44").Value = "  +2.01%  "
$ws.Range("D45").Value = "'0.00000000132"
$ws.Range("E45").Value = "  +8.32%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "2.041.70"
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("D48").Value = "'0.5222"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").Value = "'1.801"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "'9.467"
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("D51").Value = "'0.4359"
$ws.Range("E51").Value = "  +0.97%  "
